# "Add files via upload" — apply the tracked edits to Story Sharing Website.pptx
#
#  1. Slide 2 ("Idea"): merge the two runs of the last bullet into one run.
#  2. Slide 3 ("Unique Selling Point"): bump the autofit shrink + insert a new
#     bullet right after "Gap in market".
#  3. Slide 5 ("Finance Plan"): bump the autofit shrink, reword two bullets and
#     insert a new "Readers can sponsor writers" bullet before "Ad revenue".
#  4. Delete slide 4 ("Plan of website:") entirely, leaving "Finance Plan" as
#     the new slide 4.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 2 ("Idea") — merge "Writers cannot get reputation points without "
#    and "publishing stories" into a single run.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$find2 = $tr2.Find("Writers cannot get reputation points without")
$tail2 = $tr2.Characters($find2.Start, $tr2.Length - $find2.Start + 1)
$tail2.Text = "Writers cannot get reputation points without publishing stories"

# ---------------------------------------------------------------------------
# 2. Slide 3 ("Unique Selling Point")
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)

# Autofit got tighter once the new bullet was added.
try {
    $sh3.TextFrame2.AutoSize = 2
    $sh3.TextFrame2.FontScale = 70000
    $sh3.TextFrame2.LineSpaceReduction = 20000
} catch {}

$tr3 = $sh3.TextFrame.TextRange
$newText3 = "Target audience is everyone who enjoys reading and writing" + [char]13 + `
    "It is very simple to use" + [char]13 + `
    "Gap in market" + [char]13 + `
    "More accessible to amateur writers, includes a complete WYSIWYG editor so no need to learn any kind of mark up" + [char]13 + `
    "Writers can kick-start careers without resources needed for publishing" + [char]13 + `
    "Basic Design:" + [char]13
$tr3.Text = $newText3

# Restore the "Basic Design:" paragraph's no-bullet formatting.
$bd3 = $tr3.Find("Basic Design:")
$bd3.ParagraphFormat.Bullet.Visible = 0

# ---------------------------------------------------------------------------
# 3. Slide 5 ("Finance Plan") — reword two bullets and add a new one.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)

try {
    $sh5.TextFrame2.AutoSize = 2
    $sh5.TextFrame2.FontScale = 92500
    $sh5.TextFrame2.LineSpaceReduction = 20000
} catch {}

$tr5 = $sh5.TextFrame.TextRange
$newText5 = "Readers can read the first 2 pages before payment is necessary, to capture their interest" + [char]13 + `
    "Writers get 5 free stories on sign-up a we get 30% of the revenue from free stories. Writers get 1 free chapter every week to keep them coming back." + [char]13 + `
    "When writer starts paying for stories, we take 15% of revenue" + [char]13 + `
    "Writer sets price of story, but we have algorithm running showing optimal price based on story length and user reputation points" + [char]13 + `
    "Writers can subscribe which will give discounted price per story upload" + [char]13 + `
    "Readers can sponsor writers" + [char]13 + `
    "Ad revenue" + [char]13 + `
    "Ad free experience (Only £1.99 per month)"
$tr5.Text = $newText5

# ---------------------------------------------------------------------------
# 4. Remove slide 4 ("Plan of website:") — "Finance Plan" slides into its spot.
# ---------------------------------------------------------------------------
$p.Slides.Item(4).Delete()
